# Elimina EC anteriores y se agregan nuevos, se modifica base de datos
# Reorders the account-statement detail rows (B16:J45) so the data is
# grouped per worker (instead of per period), descending by period.

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Hoja1")

# Target data for rows 16..45: Doc#, Name, Period, Mora value, Base salary.
# (Type-of-doc column B and the trailing H/I/J columns are unchanged.)
$data = @(
    @("45460668", "LAIS DEL CARMEN GARCIA PEREZ", "2109", 24227, 828116),
    @("45460668", "LAIS DEL CARMEN GARCIA PEREZ", "2108", 36341, 828116),
    @("45460668", "LAIS DEL CARMEN GARCIA PEREZ", "2107", 36341, 828116),
    @("45460668", "LAIS DEL CARMEN GARCIA PEREZ", "2106", 36341, 828116),
    @("45460668", "LAIS DEL CARMEN GARCIA PEREZ", "2105", 36341, 828116),
    @("45460668", "LAIS DEL CARMEN GARCIA PEREZ", "2104", 36341, 828116),
    @("45520490", "MEREDITH MORALES CONTRERAS", "2109", 24227, 908526),
    @("45520490", "MEREDITH MORALES CONTRERAS", "2108", 36341, 908526),
    @("45520490", "MEREDITH MORALES CONTRERAS", "2107", 36341, 908526),
    @("45520490", "MEREDITH MORALES CONTRERAS", "2106", 36341, 908526),
    @("45520490", "MEREDITH MORALES CONTRERAS", "2105", 36341, 908526),
    @("45520490", "MEREDITH MORALES CONTRERAS", "2104", 36341, 908526),
    @("79951051", "JAROLD ORLANDO CASTAÑEDA", "2109", 26666, 1000000),
    @("79951051", "JAROLD ORLANDO CASTAÑEDA", "2108", 40000, 1000000),
    @("79951051", "JAROLD ORLANDO CASTAÑEDA", "2107", 40000, 1000000),
    @("79951051", "JAROLD ORLANDO CASTAÑEDA", "2106", 40000, 1000000),
    @("79951051", "JAROLD ORLANDO CASTAÑEDA", "2105", 40000, 1000000),
    @("79951051", "JAROLD ORLANDO CASTAÑEDA", "2104", 40000, 1000000),
    @("1143471376", "CHARLES RICARDO AYALA RIOS", "2109", 24227, 908526),
    @("1143471376", "CHARLES RICARDO AYALA RIOS", "2108", 36341, 908526),
    @("1143471376", "CHARLES RICARDO AYALA RIOS", "2107", 36341, 908526),
    @("1143471376", "CHARLES RICARDO AYALA RIOS", "2106", 36341, 908526),
    @("1143471376", "CHARLES RICARDO AYALA RIOS", "2105", 36341, 908526),
    @("1143471376", "CHARLES RICARDO AYALA RIOS", "2104", 36341, 908526),
    @("1070811526", "HUMBERTO LOPEZ HERAZO", "2109", 24227, 908526),
    @("1070811526", "HUMBERTO LOPEZ HERAZO", "2108", 36341, 908526),
    @("1070811526", "HUMBERTO LOPEZ HERAZO", "2107", 36341, 908526),
    @("1070811526", "HUMBERTO LOPEZ HERAZO", "2106", 36341, 908526),
    @("1070811526", "HUMBERTO LOPEZ HERAZO", "2105", 36341, 908526),
    @("1070811526", "HUMBERTO LOPEZ HERAZO", "2104", 36341, 908526)
)

$startRow = 16
for ($i = 0; $i -lt $data.Count; $i++) {
    $row = $startRow + $i
    $item = $data[$i]
    $ws.Cells.Item($row, 3).Value = $item[0]   # C: N Doc Trabajador
    $ws.Cells.Item($row, 4).Value = $item[1]   # D: Nombre Trabajador
    $ws.Cells.Item($row, 5).Value = $item[2]   # E: Periodo Mora
    $ws.Cells.Item($row, 6).Value = $item[3]   # F: Valor Mora
    $ws.Cells.Item($row, 7).Value = $item[4]   # G: Salario Basico
}
